$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Acer campestre"
$ws.Range("B3").Value = "(tree)"
$ws.Range("A4").Value = "Acer negundo"
$ws.Range("B4").Value = "(tree)"
$ws.Range("A5").Value = "Acer platanoides"
$ws.Range("B5").Value = "(tree)"
$ws.Range("A6").Value = "Acer pseudoplatanus"
$ws.Range("B6").Value = "(tree)"
$ws.Range("A7").Value = "Acer species"
$ws.Range("B7").Value = "not assigned"

$ws.Range("A11").Value = "Allium species"
$ws.Range("B11").Value = "not assigned"

$ws.Range("A17").Value = "Apiaceae species"
$ws.Range("B17").Value = "not assigned"

$ws.Range("A22").Value = "Asteraceae species"
$ws.Range("B22").Value = "not assigned"

$ws.Range("A27").Value = "Betula species"

$ws.Range("A32").Value = "Brassicaceae species"
$ws.Range("B32").Value = "not assigned"

$ws.Range("B55").Value = "(shrub/tree)"

$ws.Range("B57").Value = "(shrub/tree)"

$ws.Range("A66").Value = "Draba species"
$ws.Range("B66").Value = "not assigned"

$ws.Range("B69").Value = "(fern)"

$ws.Range("A77").Value = "Festuca ovina"
$ws.Range("A78").Value = "Festuca rubra"
$ws.Range("A79").Value = "Festuca rupicola"
$ws.Range("A80").Value = "Festuca species"

$ws.Range("B84").Value = "(tree)"

$ws.Range("A91").Value = "Geranium molle"
$ws.Range("A92").Value = "Geranium pratense"
$ws.Range("A93").Value = "Geranium pusillum"
$ws.Range("A94").Value = "Geranium pyrenaicum"
$ws.Range("A95").Value = "Geranium rotundifolium"
$ws.Range("A96").Value = "Geranium species"
$ws.Range("B96").Value = "not assigned"

$ws.Range("A123").Value = "Medicago falcata"
$ws.Range("A124").Value = "Medicago lupulina"
$ws.Range("A125").Value = "Medicago species"
$ws.Range("B125").Value = "not assigned"

$ws.Range("A145").Value = "Poaceae species"

$ws.Range("B147").Value = "(shrub/tree)"

$ws.Range("A151").Value = "Prunus avium"
$ws.Range("B151").Value = "(tree)"
$ws.Range("A152").Value = "Prunus mahaleb"
$ws.Range("B152").Value = "(shrub/tree)"
$ws.Range("A153").Value = "Prunus species"
$ws.Range("B153").Value = "not assigned"

$ws.Range("A159").Value = "Rubus caesius"
$ws.Range("B159").Value = "(shrub)"
$ws.Range("A160").Value = "Rubus idaeus"
$ws.Range("B160").Value = "(shrub)"
$ws.Range("A161").Value = "Rubus species"
$ws.Range("B161").Value = "not assigned"

$ws.Range("B164").Value = "(shrub/tree)"

$ws.Range("A167").Value = "Senecio jacobaea"
$ws.Range("B167").Value = "forb"
$ws.Range("A168").Value = "Senecio species"
$ws.Range("B168").Value = "not assigned"

$ws.Range("B192").Value = "not assigned"

$ws.Range("A202").Value = "Triticum species"

$ws.Range("A213").Value = "Vicia cracca"
$ws.Range("A214").Value = "Vicia hirsuta"
$ws.Range("A215").Value = "Vicia sativa"
$ws.Range("A216").Value = "Vicia sepium"
$ws.Range("A217").Value = "Vicia species"
$ws.Range("B217").Value = "not assigned"
